# Update "想去人数" (want-to-go count) figures in column F across the
# relevant sheets, matching the freshly generated gh-pages data snapshot.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value  = 1588
$ws.Cells.Item(3, 6).Value  = 672
$ws.Cells.Item(4, 6).Value  = 364
$ws.Cells.Item(5, 6).Value  = 5171
$ws.Cells.Item(7, 6).Value  = 10185
$ws.Cells.Item(8, 6).Value  = 263
$ws.Cells.Item(9, 6).Value  = 557
$ws.Cells.Item(11, 6).Value = 82
$ws.Cells.Item(12, 6).Value = 788

# --- Sheet "演出" (Performances) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(4, 6).Value = 15

# --- Sheet "全部类型" (All types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value  = 1588
$ws.Cells.Item(3, 6).Value  = 672
$ws.Cells.Item(4, 6).Value  = 364
$ws.Cells.Item(7, 6).Value  = 5171
$ws.Cells.Item(9, 6).Value  = 15
$ws.Cells.Item(10, 6).Value = 10185
$ws.Cells.Item(11, 6).Value = 263
$ws.Cells.Item(12, 6).Value = 557
$ws.Cells.Item(16, 6).Value = 82
$ws.Cells.Item(17, 6).Value = 788
